$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Row 8: new "log(cpc-2)" label + formula ---
$ws.Range("A8").Value = "log(cpc-2)"
$ws.Range("C8").Formula = "=LOG(C7,2)"

# Pre-create the "applyNumberFormat" style (cellXfs index 5) ahead of the
# quotePrefix style (index 6) below, matching the authoring order in the
# target file (I37 was number-formatted before B24/B26 got their leading
# apostrophe).
$ws.Range("I37").Value = 1024
$ws.Range("I37").NumberFormat = "General"

# --- Row 21-22: headers for the new resource-utilization matrix ---
$ws.Range("B21").Value = "Other Modules"
$ws.Range("D21").Value = "Multipliers"
$ws.Range("E21").Value = "Add/Comp"
$ws.Range("F21").Value = "Add/Comp"
$ws.Range("G21").Value = "DFFs"
$ws.Range("H21").Value = "MUXes"
$ws.Range("I21").Value = "Mem cells"
$ws.Range("J21").Value = "Mem each cell"
$ws.Range("K21").Value = "Total mem bits"
$ws.Range("N21").Value = "Comments"

$ws.Range("D22").Value = "width bits"
$ws.Range("E22").Value = "width bits"
$ws.Range("F22").Value = "1 bit"
$ws.Range("G22").Value = "1 bit"
$ws.Range("H22").Value = "N to 1"

# --- Row 23: Multiplier ---
$ws.Range("A23").Value = "Multiplier"
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("F23").Formula = "=1+I2+1+1+I2+1"

# --- Row 24: Mult Set ---
$ws.Range("A24").Value = "Mult Set"
$ws.Range("B24").Value = "'" + '"z" Mulitpliers'

# --- Row 25: Adder ---
$ws.Range("A25").Value = "Adder"
$ws.Range("E25").Value = 1
$ws.Range("F25").Formula = "=1+1+1"

# --- Row 26: costterm set ---
$ws.Range("A26").Value = "costterm set"
$ws.Range("B26").Value = "'" + '"z" Adders'

# --- Row 27: counter ---
$ws.Range("A27").Value = "counter"
$ws.Range("F27").Formula = "=4*2"
$ws.Range("G27").Formula = "=F27/2"
$ws.Range("N27").Value = "Assume max is 4b"

# --- Row 28: cycleblockcnt ---
$ws.Range("A28").Value = "cycleblockcnt"
$ws.Range("F28").Formula = "=C8*2"
$ws.Range("G28").Formula = "=C8"

# --- Row 29: mux ---
$ws.Range("A29").Value = "mux"
$ws.Range("H29").Formula = "=H2"

# --- Row 30: muxset ---
$ws.Range("A30").Value = "muxset"
$ws.Range("B30").Value = "M muxes"

# --- Row 31: maxfinder ---
$ws.Range("A31").Value = "maxfinder"
$ws.Range("E31").Value = 2

# --- Row 32: maxfinderset ---
$ws.Range("A32").Value = "maxfinderset"
$ws.Range("B32").Value = "2N maxfinders"
$ws.Range("F32").Value = 1

# --- Row 33: DFF ---
$ws.Range("A33").Value = "DFF"
$ws.Range("G33").Value = "width"

# --- Row 34: shiftreg ---
$ws.Range("A34").Value = "shiftreg"
$ws.Range("G34").Value = "width*depth"

# --- Row 36: sigmoid table ---
$ws.Range("A36").Value = "sigmoid table"
$ws.Range("E36").Value = 2
$ws.Range("I36").Value = 1024
$ws.Range("J36").Formula = "=J2"
$ws.Range("K36").Formula = "=I36*J36"

# --- Row 37: sp table (I37 value/format were set above) ---
$ws.Range("A37").Value = "sp table"
$ws.Range("E37").Value = 1
$ws.Range("J37").Formula = "=J2-2"
$ws.Range("K37").Formula = "=I37*J37"

# --- Row 39: interleaverset ---
$ws.Range("A39").Value = "interleaverset"
$ws.Range("G39").Formula = "=C3*C6*LOG(B2/C6,2)"

# --- Row 41: sig function ---
$ws.Range("A41").Value = "sig function"
$ws.Range("B41").Value = "2*fi+1 Adders (width_TA bits)"
$ws.Range("F41").Value = 2

# --- Row 42-43 ---
$ws.Range("B42").Value = "1 sigmoid table"
$ws.Range("B43").Value = "1 sp table"

# --- Row 45-46: FF proc ---
$ws.Range("A45").Value = "FF proc"
$ws.Range("B45").Value = "z Multipliers"
$ws.Range("B46").Value = "z/fi sig functions"

# --- Row 48-49: BP proc ---
$ws.Range("A48").Value = "BP proc"
$ws.Range("B48").Value = "2z Multipliers"
$ws.Range("B49").Value = "z Adders"

# --- Row 51: UP proc ---
$ws.Range("A51").Value = "UP proc"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 12.33203125
$ws.Columns.Item(3).ColumnWidth = 10.83203125
$ws.Columns.Item(8).ColumnWidth = 12.83203125
$ws.Columns.Item(9).ColumnWidth = 13.33203125

# --- View: zoom + selection ---
$excel.ActiveWindow.Zoom = 125
$ws.Range("G48").Select()
